# Automatic update of files.
# Rows 2-7 of the sheet get their A/B/E/F/G/H/Q/R column values rotated
# among each other (a single 6-cycle: 2 -> 3 -> 5 -> 7 -> 6 -> 4 -> 2).
# We compute the new values first (from the *current* cell contents) and
# then write them all back, so the operation behaves like a simultaneous
# row swap rather than a sequence of overwrites that clobber source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# after_row -> before_row (i.e. new content at $after comes from old content at $before)
$mapping = @{
    2 = 4
    3 = 2
    4 = 6
    5 = 3
    6 = 7
    7 = 5
}

# Snapshot the current ("before") values for every relevant cell first.
$snapshot = @{}
foreach ($row in $mapping.Values) {
    if (-not $snapshot.ContainsKey($row)) {
        $rowData = @{}
        foreach ($col in $cols) {
            $rowData[$col] = $ws.Range("$col$row").Value()
        }
        $snapshot[$row] = $rowData
    }
}

# Now write the rotated values into their destination rows.
foreach ($afterRow in $mapping.Keys) {
    $beforeRow = $mapping[$afterRow]
    $rowData = $snapshot[$beforeRow]
    foreach ($col in $cols) {
        $ws.Range("$col$afterRow").Value = $rowData[$col]
    }
}
